$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column name for the "Round" table's DateTime column to "PlayTime"
$ws.Range("B4").Value = "PlayTime"

# Reflect the active selection recorded at save time
$ws.Range("B4").Select()
